$d = $word.ActiveDocument

# Reusable pkg:package wrapper needed by Range.InsertXML
$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Locate the existing "Select * from student;" paragraph by its text so the
# script does not depend on a fixed paragraph index.
function Find-StudentParagraph($doc) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like "Select * from student;*") {
            return $p
        }
    }
    return $null
}

# --- Insert a new empty paragraph BEFORE "Select * from student;" ---
$studentPara = Find-StudentParagraph $d
$insertStart = $d.Range($studentPara.Range.Start, $studentPara.Range.Start)
$emptyParaXml = $pkgHeader + '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' + $pkgFooter
$insertStart.InsertXML($emptyParaXml)

# Re-fetch the paragraph: the previous reference/range does not track the
# document after the mutation above, so look it up again by its text.
$studentPara = Find-StudentParagraph $d

# --- Insert a new paragraph with text AFTER "Select * from student;" ---
$insertEnd = $d.Range($studentPara.Range.End, $studentPara.Range.End)
$newParaXml = $pkgHeader + '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Select * from profeesor;</w:t></w:r></w:p>' + $pkgFooter
$insertEnd.InsertXML($newParaXml)

Write-Output "Paragraph count: $($d.Paragraphs.Count)"
foreach ($p in $d.Paragraphs) {
    Write-Output "[$($p.Range.Text)]"
}
